$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-decimal Price cells to remain Text (matches source data which stores
# prices/volumes as literal strings, not numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "60.195.00"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "2.617.90"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "568.82"
$ws.Range("E5").Value = "  +5.57%  "
$ws.Range("D6").Value = "145.55"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +4.27%  "
$ws.Range("D9").Value = "2.632.38"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("E12").Value = "  +9.19%  "
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("D14").Value = "3.077.82"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "60.162.48"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "22.08"
$ws.Range("E16").Value = "  +6.29%  "
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "2.629.87"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "4.52"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "340.02"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "10.38"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "65.70"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").Value = "0.448"
$ws.Range("E25").Value = "  +6.51%  "
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("D29").Value = "0.0₃0793"
$ws.Range("E29").Value = "  +7.46%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "6.13"
$ws.Range("E32").Value = "  +2.95%  "
$ws.Range("D33").Value = "160.61"
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").Value = "0.881"
$ws.Range("E37").Value = "  +4.23%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "0.885"
$ws.Range("E38").Value = "  +7.45%  "
$ws.Range("D39").Value = "37.50"
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("E40").Value = "  +5.99%  "
$ws.Range("D41").Value = "296.83"
$ws.Range("E41").Value = "  +4.74%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "0.0980"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("D46").Value = "0.0540"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "19.29"
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").Value = "10.67"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").Value = "126.71"
$ws.Range("E49").Value = "  +14.89%  "
$ws.Range("D50").Value = "0.0234"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").Value = "1.958.65"
$ws.Range("E51").Value = "  +1.99%  "
